$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns now report "handed back" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: status updated, handback datetime refreshed, error cleared ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-10-18 11:47:38"
$zhcn.Range("P2").Value = ""

# --- de-de sheet: status updated, handback datetime refreshed, error cleared ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-10-18 11:47:55"
$dede.Range("P2").Value = ""

# --- Column widths: Status columns widen to fit the longer message, the
#     now-empty Error Detail columns narrow back down. ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
